$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2,2).Value = '2024-06-23'
$ws1.Cells.Item(2,3).Value = '赣州·清风霁月·光夜only'
$ws1.Cells.Item(2,4).Value = '平安大道 麋鹿LiveHouse'
$ws1.Cells.Item(2,5).Value = '2024.06.23 14:00-06.23 20:00'
$ws1.Cells.Item(2,6).Value = 91
$ws1.Cells.Item(2,7).Value = '不可售'
$ws1.Cells.Item(2,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86993'
$ws1.Cells.Item(2,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/PklWR8EP1717429316070.jpeg'

$ws1.Cells.Item(3,2).Value = '2024-06-29'
$ws1.Cells.Item(3,3).Value = '南昌·第五人格only'
$ws1.Cells.Item(3,4).Value = '高处见美好生活公园 百家喜宴高新店'
$ws1.Cells.Item(3,5).Value = '2024.06.29 10:00-06.29 17:00'
$ws1.Cells.Item(3,6).Value = 325
$ws1.Cells.Item(3,7).Value = 65
$ws1.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87043'
$ws1.Cells.Item(3,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/zir2PYz81717071721569.jpeg'

$ws1.Cells.Item(4,2).Value = '2024-06-29'
$ws1.Cells.Item(4,3).Value = '萍乡·BM次元盛典运动番only'
$ws1.Cells.Item(4,4).Value = '康庄路3号 萍乡梅园国际大酒店'
$ws1.Cells.Item(4,5).Value = '2024.06.29 10:00-06.29 17:00'
$ws1.Cells.Item(4,6).Value = 278
$ws1.Cells.Item(4,7).Value = 55
$ws1.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85192'
$ws1.Cells.Item(4,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png'

$ws1.Cells.Item(5,2).Value = '2024-06-30'
$ws1.Cells.Item(5,3).Value = '南昌·ChinastyleCOSPLAY  '
$ws1.Cells.Item(5,4).Value = '真君路888号 南昌华侨城玩美公园'
$ws1.Cells.Item(5,5).Value = '2024.06.30 09:30-07.02 17:30'
$ws1.Cells.Item(5,6).Value = 1193
$ws1.Cells.Item(5,7).Value = 65
$ws1.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87045'
$ws1.Cells.Item(5,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/wajWy7ID1717149642528.jpeg'

$ws1.Cells.Item(6,2).Value = '2024-06-30'
$ws1.Cells.Item(6,3).Value = '宜春·BM次元盛典运动番only'
$ws1.Cells.Item(6,4).Value = '鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)'
$ws1.Cells.Item(6,5).Value = '2024.06.30 10:00-06.30 17:00'
$ws1.Cells.Item(6,6).Value = 271
$ws1.Cells.Item(6,7).Value = 55
$ws1.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84636'
$ws1.Cells.Item(6,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png'

$ws1.Cells.Item(7,2).Value = '2024-07-06'
$ws1.Cells.Item(7,3).Value = '南昌·次元星球动漫游戏展'
$ws1.Cells.Item(7,4).Value = '龙蟠街666号融创茂1层 融创茂'
$ws1.Cells.Item(7,5).Value = '2024.07.06 10:00-07.06 17:00'
$ws1.Cells.Item(7,6).Value = 25
$ws1.Cells.Item(7,7).Value = '不可售'
$ws1.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86405'
$ws1.Cells.Item(7,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg'

$ws1.Cells.Item(8,2).Value = '2024-07-06'
$ws1.Cells.Item(8,3).Value = '鹰潭·BM次元盛典运动番only'
$ws1.Cells.Item(8,4).Value = '体育馆东路2号九小隔壁 忆江南•宴会楼'
$ws1.Cells.Item(8,5).Value = '2024.07.06 10:00-07.06 17:00'
$ws1.Cells.Item(8,6).Value = 58
$ws1.Cells.Item(8,7).Value = 55
$ws1.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85997'
$ws1.Cells.Item(8,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png'

$ws1.Cells.Item(9,2).Value = '2024-07-07'
$ws1.Cells.Item(9,3).Value = '赣州·BM次元盛典运动番only'
$ws1.Cells.Item(9,4).Value = '米瑞金路2口0号上客天下1楼 上客天下.老虔州'
$ws1.Cells.Item(9,5).Value = '2024.07.07 10:00-07.07 17:00'
$ws1.Cells.Item(9,6).Value = 40
$ws1.Cells.Item(9,7).Value = 55
$ws1.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86602'
$ws1.Cells.Item(9,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png'

$ws1.Cells.Item(10,2).Value = '2024-07-12'
$ws1.Cells.Item(10,3).Value = '新余·2024第三届MG动漫嘉年华'
$ws1.Cells.Item(10,4).Value = '仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅'
$ws1.Cells.Item(10,5).Value = '2024.07.12 10:00-07.13 17:30'
$ws1.Cells.Item(10,6).Value = 134
$ws1.Cells.Item(10,7).Value = 55
$ws1.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86536'
$ws1.Cells.Item(10,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg'

$ws1.Cells.Item(11,2).Value = '2024-07-13'
$ws1.Cells.Item(11,3).Value = '南昌·SuperComic动漫游戏博览会'
$ws1.Cells.Item(11,4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Cells.Item(11,5).Value = '2024.07.13 09:00-07.14 17:00'
$ws1.Cells.Item(11,6).Value = 3388
$ws1.Cells.Item(11,7).Value = 65
$ws1.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86992'
$ws1.Cells.Item(11,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/wQTAjelJ1717642148929.jpeg'

$ws1.Cells.Item(12,2).Value = '2024-07-13'
$ws1.Cells.Item(12,3).Value = '南昌·SuperComic配音演员刘明月专场见面会'
$ws1.Cells.Item(12,4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Cells.Item(12,5).Value = '2024.07.13 09:00-07.13 17:00'
$ws1.Cells.Item(12,6).Value = 110
$ws1.Cells.Item(12,7).Value = 168
$ws1.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87570'
$ws1.Cells.Item(12,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/1D1reIl81718609013880.png'

$ws1.Cells.Item(13,2).Value = '2024-07-13'
$ws1.Cells.Item(13,3).Value = '南昌·THO-梦违赣鄱荟萃·叁~幻想Strawberry~!!'
$ws1.Cells.Item(13,4).Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws1.Cells.Item(13,5).Value = '2024.07.13 09:30-07.13 17:30'
$ws1.Cells.Item(13,6).Value = 77
$ws1.Cells.Item(13,7).Value = 65
$ws1.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87668'
$ws1.Cells.Item(13,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/Bk9cYryT1718360290362.jpeg'

$ws1.Cells.Item(14,2).Value = '2024-07-13'
$ws1.Cells.Item(14,3).Value = '宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华'
$ws1.Cells.Item(14,4).Value = '宜春国际商贸城会展中心 宜春国际商贸城会展中心'
$ws1.Cells.Item(14,5).Value = '2024.07.13 10:00-07.14 17:00'
$ws1.Cells.Item(14,6).Value = 63
$ws1.Cells.Item(14,7).Value = 55
$ws1.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86667'
$ws1.Cells.Item(14,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg'

$ws1.Cells.Item(15,2).Value = '2024-07-13'
$ws1.Cells.Item(15,3).Value = '赣州·十万伏特-次元交流会（夏）'
$ws1.Cells.Item(15,4).Value = '梅关大道36-16号 麋鹿星球艺术中心'
$ws1.Cells.Item(15,5).Value = '2024.07.13 09:30-07.13 17:00'
$ws1.Cells.Item(15,6).Value = 36
$ws1.Cells.Item(15,7).Value = 45
$ws1.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87597'
$ws1.Cells.Item(15,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/87yQ4Hmf1718681348727.jpeg'

$ws1.Cells.Item(16,2).Value = '2024-07-14'
$ws1.Cells.Item(16,3).Value = '南昌·赛马娘ONLY'
$ws1.Cells.Item(16,4).Value = '洪城路99号 锦都皇冠酒店(八一广场火车站店)'
$ws1.Cells.Item(16,5).Value = '2024.07.14 09:00-07.14 17:30'
$ws1.Cells.Item(16,6).Value = 50
$ws1.Cells.Item(16,7).Value = 68
$ws1.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87367'
$ws1.Cells.Item(16,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/wXQuIKtu1718165450704.png'

$ws1.Cells.Item(17,2).Value = '2024-07-14'
$ws1.Cells.Item(17,3).Value = '吉安·COMIC LIFE次元假日05'
$ws1.Cells.Item(17,4).Value = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws1.Cells.Item(17,5).Value = '2024.07.14 09:00-07.14 18:00'
$ws1.Cells.Item(17,6).Value = 577
$ws1.Cells.Item(17,7).Value = 52.1
$ws1.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85924'
$ws1.Cells.Item(17,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg'

$ws1.Cells.Item(18,2).Value = '2024-07-19'
$ws1.Cells.Item(18,3).Value = '九江·第一届Loading加载中动漫展'
$ws1.Cells.Item(18,4).Value = '湓浦街道大中路339号 百嘉洲际酒店'
$ws1.Cells.Item(18,5).Value = '2024.07.19 09:00-07.21 17:00'
$ws1.Cells.Item(18,6).Value = 62
$ws1.Cells.Item(18,7).Value = 36.6
$ws1.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87787'
$ws1.Cells.Item(18,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/TH3lVD5G1718158901239.jpeg'

$ws1.Cells.Item(19,2).Value = '2024-07-19'
$ws1.Cells.Item(19,3).Value = '赣州·第四届赣州半夏动漫展'
$ws1.Cells.Item(19,4).Value = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws1.Cells.Item(19,5).Value = '2024.07.19 10:00-07.21 17:00'
$ws1.Cells.Item(19,6).Value = 686
$ws1.Cells.Item(19,7).Value = 55
$ws1.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86587'
$ws1.Cells.Item(19,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg'

$ws1.Cells.Item(20,2).Value = '2024-07-20'
$ws1.Cells.Item(20,3).Value = '南昌·漫拥动漫嘉年华Pro-追光启航'
$ws1.Cells.Item(20,4).Value = '小蓝南路420号 洪州体育馆'
$ws1.Cells.Item(20,5).Value = '2024.07.20 09:00-07.21 17:00'
$ws1.Cells.Item(20,6).Value = 198
$ws1.Cells.Item(20,7).Value = 52.5
$ws1.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85796'
$ws1.Cells.Item(20,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png'

$ws1.Cells.Item(21,2).Value = '2024-07-21'
$ws1.Cells.Item(21,3).Value = '乐平·CY境界次元动漫夏时庆'
$ws1.Cells.Item(21,4).Value = '翥山西路182号 佳佳基大酒店'
$ws1.Cells.Item(21,5).Value = '2024.07.21 10:00-07.21 17:00'
$ws1.Cells.Item(21,6).Value = 111
$ws1.Cells.Item(21,7).Value = 35
$ws1.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86768'
$ws1.Cells.Item(21,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png'

$ws1.Cells.Item(22,2).Value = '2024-07-21'
$ws1.Cells.Item(22,3).Value = '九江·SXD动漫嘉年华'
$ws1.Cells.Item(22,4).Value = '湓浦街道大中路339号 百嘉洲际酒店'
$ws1.Cells.Item(22,5).Value = '2024.07.21 10:00-07.21 17:30'
$ws1.Cells.Item(22,6).Value = 54
$ws1.Cells.Item(22,7).Value = 45
$ws1.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86832'
$ws1.Cells.Item(22,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg'

$ws1.Cells.Item(23,2).Value = '2024-07-21'
$ws1.Cells.Item(23,3).Value = '抚州·临次元08·盛夏动漫狂欢节'
$ws1.Cells.Item(23,4).Value = '伍塘路1098号 乐课篮球公园'
$ws1.Cells.Item(23,5).Value = '2024.07.21 10:00-07.21 16:00'
$ws1.Cells.Item(23,6).Value = 50
$ws1.Cells.Item(23,7).Value = 50
$ws1.Cells.Item(23,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87763'
$ws1.Cells.Item(23,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/6qgetbCh1718720523395.jpeg'

$ws1.Cells.Item(24,2).Value = '2024-07-21'
$ws1.Cells.Item(24,3).Value = '萍乡·NL14动漫游戏展·夏日狂想曲'
$ws1.Cells.Item(24,4).Value = '公园南路168号(近工行城北分理处) 梅生嘉华酒店'
$ws1.Cells.Item(24,5).Value = '2024.07.21 10:00-07.21 17:00'
$ws1.Cells.Item(24,6).Value = 59
$ws1.Cells.Item(24,7).Value = 40
$ws1.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86658'
$ws1.Cells.Item(24,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg'

$ws1.Cells.Item(25,2).Value = '2024-07-26'
$ws1.Cells.Item(25,3).Value = '南昌·萌卡动漫展'
$ws1.Cells.Item(25,4).Value = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$ws1.Cells.Item(25,5).Value = '2024.07.26 09:00-07.28 17:00'
$ws1.Cells.Item(25,6).Value = 2447
$ws1.Cells.Item(25,7).Value = 58.5
$ws1.Cells.Item(25,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86776'
$ws1.Cells.Item(25,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg'

$ws1.Cells.Item(26,2).Value = '2024-07-27'
$ws1.Cells.Item(26,3).Value = '江西·次元星河动漫游戏嘉年华'
$ws1.Cells.Item(26,4).Value = '九龙大道1177号 南昌绿地国际博览中心'
$ws1.Cells.Item(26,5).Value = '2024.07.27 10:00-07.28 17:00'
$ws1.Cells.Item(26,6).Value = 4969
$ws1.Cells.Item(26,7).Value = 69
$ws1.Cells.Item(26,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85493'
$ws1.Cells.Item(26,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png'

$ws1.Cells.Item(27,2).Value = '2024-07-27'
$ws1.Cells.Item(27,3).Value = '赣州·马娘only'
$ws1.Cells.Item(27,4).Value = '火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)'
$ws1.Cells.Item(27,5).Value = '2024.07.27 09:00-07.27 17:00'
$ws1.Cells.Item(27,6).Value = 32
$ws1.Cells.Item(27,7).Value = 60
$ws1.Cells.Item(27,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86772'
$ws1.Cells.Item(27,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png'

$ws1.Cells.Item(28,2).Value = '2024-07-28'
$ws1.Cells.Item(28,3).Value = '赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会'
$ws1.Cells.Item(28,4).Value = '兴国路恒大帝景西门 江西长庚控股有限公司'
$ws1.Cells.Item(28,5).Value = '2024.07.28 11:00-07.28 17:00'
$ws1.Cells.Item(28,6).Value = 69
$ws1.Cells.Item(28,7).Value = 56
$ws1.Cells.Item(28,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85688'
$ws1.Cells.Item(28,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png'

$ws1.Cells.Item(29,2).Value = '2024-07-30'
$ws1.Cells.Item(29,3).Value = '宜春·第三十五届静卿国风动漫文化展览会'
$ws1.Cells.Item(29,4).Value = '宜阳大道19号(交通银行旁) 宜春安缦文华酒店'
$ws1.Cells.Item(29,5).Value = '2024.07.30 09:00-07.30 17:00'
$ws1.Cells.Item(29,6).Value = 472
$ws1.Cells.Item(29,7).Value = 45
$ws1.Cells.Item(29,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86684'
$ws1.Cells.Item(29,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg'

$ws1.Cells.Item(30,2).Value = '2024-08-03'
$ws1.Cells.Item(30,3).Value = '南昌·幻梦境国际动漫游戏嘉年华1th'
$ws1.Cells.Item(30,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws1.Cells.Item(30,5).Value = '2024.08.03 09:00-08.04 17:30'
$ws1.Cells.Item(30,6).Value = 1277
$ws1.Cells.Item(30,7).Value = 64
$ws1.Cells.Item(30,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83980'
$ws1.Cells.Item(30,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'

$ws1.Cells.Item(31,2).Value = '2024-08-03'
$ws1.Cells.Item(31,3).Value = '吉安·COMIC LIFE周年庆典'
$ws1.Cells.Item(31,4).Value = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws1.Cells.Item(31,5).Value = '2024.08.03 09:30-08.03 18:00'
$ws1.Cells.Item(31,6).Value = 277
$ws1.Cells.Item(31,7).Value = 46.6
$ws1.Cells.Item(31,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87164'
$ws1.Cells.Item(31,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg'

$ws1.Cells.Item(32,2).Value = '2024-08-03'
$ws1.Cells.Item(32,3).Value = '景德镇·第十五届瓷都ACG动漫游戏博览会'
$ws1.Cells.Item(32,4).Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws1.Cells.Item(32,5).Value = '2024.08.03 09:00-08.04 17:00'
$ws1.Cells.Item(32,6).Value = 2199
$ws1.Cells.Item(32,7).Value = 55
$ws1.Cells.Item(32,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86341'
$ws1.Cells.Item(32,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png'

$ws1.Cells.Item(33,2).Value = '2024-08-03'
$ws1.Cells.Item(33,3).Value = '景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票'
$ws1.Cells.Item(33,4).Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws1.Cells.Item(33,5).Value = '2024.08.03 08:30-08.03 17:00'
$ws1.Cells.Item(33,6).Value = 570
$ws1.Cells.Item(33,7).Value = '已售罄'
$ws1.Cells.Item(33,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85981'
$ws1.Cells.Item(33,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png'

$ws1.Cells.Item(34,2).Value = '2024-08-03'
$ws1.Cells.Item(34,3).Value = '樟树·第二届静卿国风动漫文化展览会'
$ws1.Cells.Item(34,4).Value = '杏佛路89号 樟树银河国际酒店'
$ws1.Cells.Item(34,5).Value = '2024.08.03 09:00-08.03 17:00'
$ws1.Cells.Item(34,6).Value = 484
$ws1.Cells.Item(34,7).Value = 45
$ws1.Cells.Item(34,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86683'
$ws1.Cells.Item(34,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg'

$ws1.Cells.Item(35,2).Value = '2024-08-03'
$ws1.Cells.Item(35,3).Value = '萍乡·AU9夏至国漫展'
$ws1.Cells.Item(35,4).Value = '安源中大道17号 壹号公馆（萍乡）'
$ws1.Cells.Item(35,5).Value = '2024.08.03 10:00-08.03 17:00'
$ws1.Cells.Item(35,6).Value = 75
$ws1.Cells.Item(35,7).Value = 45
$ws1.Cells.Item(35,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86453'
$ws1.Cells.Item(35,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/hm1EACno1718936156944.jpeg'

$ws1.Cells.Item(36,2).Value = '2024-08-03'
$ws1.Cells.Item(36,3).Value = '赣州·第一届环梦动漫游戏嘉年华'
$ws1.Cells.Item(36,4).Value = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws1.Cells.Item(36,5).Value = '2024.08.03 09:00-08.05 17:00'
$ws1.Cells.Item(36,6).Value = 85
$ws1.Cells.Item(36,7).Value = 36.6
$ws1.Cells.Item(36,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87449'
$ws1.Cells.Item(36,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/iC3PuUfR1717740188790.jpeg'

$ws1.Cells.Item(37,2).Value = '2024-08-04'
$ws1.Cells.Item(37,3).Value = '上饶·第十五届IX Group国风嘉年华暨十周年庆典'
$ws1.Cells.Item(37,4).Value = '高铁经济试验区凤凰东大道16号7幢 上饶饶商金茂诚悦酒店(上饶高铁站)'
$ws1.Cells.Item(37,5).Value = '2024.08.04 09:30-08.04 17:30'
$ws1.Cells.Item(37,6).Value = 156
$ws1.Cells.Item(37,7).Value = 60
$ws1.Cells.Item(37,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87225'
$ws1.Cells.Item(37,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/l5fIXZSX1717562269098.jpeg'

$ws1.Cells.Item(38,2).Value = '2024-08-04'
$ws1.Cells.Item(38,3).Value = '九江·第一届异次元动漫嘉年华'
$ws1.Cells.Item(38,4).Value = '长虹西大道兴城广场99号 九江半岛宾馆'
$ws1.Cells.Item(38,5).Value = '2024.08.04 08:00-08.04 17:00'
$ws1.Cells.Item(38,6).Value = 305
$ws1.Cells.Item(38,7).Value = 45
$ws1.Cells.Item(38,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84407'
$ws1.Cells.Item(38,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg'

$ws1.Cells.Item(39,2).Value = '2024-08-06'
$ws1.Cells.Item(39,3).Value = '南昌·第一届异次元动漫嘉年华'
$ws1.Cells.Item(39,4).Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws1.Cells.Item(39,5).Value = '2024.08.06 08:00-08.06 17:00'
$ws1.Cells.Item(39,6).Value = 452
$ws1.Cells.Item(39,7).Value = 55
$ws1.Cells.Item(39,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84102'
$ws1.Cells.Item(39,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg'

$ws1.Cells.Item(40,2).Value = '2024-08-08'
$ws1.Cells.Item(40,3).Value = '赣州·第二届异次元动漫嘉年华'
$ws1.Cells.Item(40,4).Value = '金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆'
$ws1.Cells.Item(40,5).Value = '2024.08.08 08:00-08.08 17:00'
$ws1.Cells.Item(40,6).Value = 777
$ws1.Cells.Item(40,7).Value = 45
$ws1.Cells.Item(40,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84184'
$ws1.Cells.Item(40,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg'

$ws1.Cells.Item(41,2).Value = '2024-08-10'
$ws1.Cells.Item(41,3).Value = '南昌·花绒万兽第二聚'
$ws1.Cells.Item(41,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws1.Cells.Item(41,5).Value = '2024.08.10 10:00-08.11 17:00'
$ws1.Cells.Item(41,6).Value = 25
$ws1.Cells.Item(41,7).Value = 188
$ws1.Cells.Item(41,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87600'
$ws1.Cells.Item(41,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/i0Ojsne01718693886054.png'

$ws1.Cells.Item(42,2).Value = '2024-08-10'
$ws1.Cells.Item(42,3).Value = '高安·第二届静卿国风动漫文化展览会'
$ws1.Cells.Item(42,4).Value = '华林中路606号 高安华鼎国际大酒店'
$ws1.Cells.Item(42,5).Value = '2024.08.10 09:00-08.10 17:00'
$ws1.Cells.Item(42,6).Value = 449
$ws1.Cells.Item(42,7).Value = 45
$ws1.Cells.Item(42,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86682'
$ws1.Cells.Item(42,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg'

$ws1.Cells.Item(43,2).Value = '2024-08-15'
$ws1.Cells.Item(43,3).Value = '上饶·次元重现夏日嘉年华'
$ws1.Cells.Item(43,4).Value = '普济巷地委大院北侧约90米 四季体育运动馆'
$ws1.Cells.Item(43,5).Value = '2024.08.15 09:30-08.15 17:30'
$ws1.Cells.Item(43,6).Value = 27
$ws1.Cells.Item(43,7).Value = 48
$ws1.Cells.Item(43,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87679'
$ws1.Cells.Item(43,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/fxlKV2SL1718784421064.jpeg'

$ws1.Cells.Item(44,2).Value = '2024-08-24'
$ws1.Cells.Item(44,3).Value = '南昌·第四届龙年动漫展——暑假最后的狂欢'
$ws1.Cells.Item(44,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws1.Cells.Item(44,5).Value = '2024.08.24 10:00-08.25 18:00'
$ws1.Cells.Item(44,6).Value = 458
$ws1.Cells.Item(44,7).Value = 45
$ws1.Cells.Item(44,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87135'
$ws1.Cells.Item(44,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'

$ws1.Rows.Item(45).Delete()

$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2,2).Value = '2024-06-23'
$ws4.Cells.Item(2,3).Value = '赣州·清风霁月·光夜only'
$ws4.Cells.Item(2,4).Value = '平安大道 麋鹿LiveHouse'
$ws4.Cells.Item(2,5).Value = '2024.06.23 14:00-06.23 20:00'
$ws4.Cells.Item(2,6).Value = 91
$ws4.Cells.Item(2,7).Value = '不可售'
$ws4.Cells.Item(2,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86993'
$ws4.Cells.Item(2,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/PklWR8EP1717429316070.jpeg'

$ws4.Cells.Item(3,2).Value = '2024-06-29'
$ws4.Cells.Item(3,3).Value = '南昌·第五人格only'
$ws4.Cells.Item(3,4).Value = '高处见美好生活公园 百家喜宴高新店'
$ws4.Cells.Item(3,5).Value = '2024.06.29 10:00-06.29 17:00'
$ws4.Cells.Item(3,6).Value = 325
$ws4.Cells.Item(3,7).Value = 65
$ws4.Cells.Item(3,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87043'
$ws4.Cells.Item(3,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/zir2PYz81717071721569.jpeg'

$ws4.Cells.Item(4,2).Value = '2024-06-29'
$ws4.Cells.Item(4,3).Value = '萍乡·BM次元盛典运动番only'
$ws4.Cells.Item(4,4).Value = '康庄路3号 萍乡梅园国际大酒店'
$ws4.Cells.Item(4,5).Value = '2024.06.29 10:00-06.29 17:00'
$ws4.Cells.Item(4,6).Value = 278
$ws4.Cells.Item(4,7).Value = 55
$ws4.Cells.Item(4,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85192'
$ws4.Cells.Item(4,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/byoupYK21714294780383.png'

$ws4.Cells.Item(5,2).Value = '2024-06-30'
$ws4.Cells.Item(5,3).Value = '南昌·ChinastyleCOSPLAY  '
$ws4.Cells.Item(5,4).Value = '真君路888号 南昌华侨城玩美公园'
$ws4.Cells.Item(5,5).Value = '2024.06.30 09:30-07.02 17:30'
$ws4.Cells.Item(5,6).Value = 1193
$ws4.Cells.Item(5,7).Value = 65
$ws4.Cells.Item(5,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87045'
$ws4.Cells.Item(5,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/wajWy7ID1717149642528.jpeg'

$ws4.Cells.Item(6,2).Value = '2024-06-30'
$ws4.Cells.Item(6,3).Value = '宜春·BM次元盛典运动番only'
$ws4.Cells.Item(6,4).Value = '鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)'
$ws4.Cells.Item(6,5).Value = '2024.06.30 10:00-06.30 17:00'
$ws4.Cells.Item(6,6).Value = 271
$ws4.Cells.Item(6,7).Value = 55
$ws4.Cells.Item(6,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84636'
$ws4.Cells.Item(6,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png'

$ws4.Cells.Item(7,2).Value = '2024-07-06'
$ws4.Cells.Item(7,3).Value = '南昌·次元星球动漫游戏展'
$ws4.Cells.Item(7,4).Value = '龙蟠街666号融创茂1层 融创茂'
$ws4.Cells.Item(7,5).Value = '2024.07.06 10:00-07.06 17:00'
$ws4.Cells.Item(7,6).Value = 25
$ws4.Cells.Item(7,7).Value = '不可售'
$ws4.Cells.Item(7,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86405'
$ws4.Cells.Item(7,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg'

$ws4.Cells.Item(8,2).Value = '2024-07-06'
$ws4.Cells.Item(8,3).Value = '鹰潭·BM次元盛典运动番only'
$ws4.Cells.Item(8,4).Value = '体育馆东路2号九小隔壁 忆江南•宴会楼'
$ws4.Cells.Item(8,5).Value = '2024.07.06 10:00-07.06 17:00'
$ws4.Cells.Item(8,6).Value = 58
$ws4.Cells.Item(8,7).Value = 55
$ws4.Cells.Item(8,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85997'
$ws4.Cells.Item(8,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png'

$ws4.Cells.Item(9,2).Value = '2024-07-07'
$ws4.Cells.Item(9,3).Value = '赣州·BM次元盛典运动番only'
$ws4.Cells.Item(9,4).Value = '米瑞金路2口0号上客天下1楼 上客天下.老虔州'
$ws4.Cells.Item(9,5).Value = '2024.07.07 10:00-07.07 17:00'
$ws4.Cells.Item(9,6).Value = 40
$ws4.Cells.Item(9,7).Value = 55
$ws4.Cells.Item(9,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86602'
$ws4.Cells.Item(9,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png'

$ws4.Cells.Item(10,2).Value = '2024-07-12'
$ws4.Cells.Item(10,3).Value = '新余·2024第三届MG动漫嘉年华'
$ws4.Cells.Item(10,4).Value = '仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅'
$ws4.Cells.Item(10,5).Value = '2024.07.12 10:00-07.13 17:30'
$ws4.Cells.Item(10,6).Value = 134
$ws4.Cells.Item(10,7).Value = 55
$ws4.Cells.Item(10,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86536'
$ws4.Cells.Item(10,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg'

$ws4.Cells.Item(11,2).Value = '2024-07-13'
$ws4.Cells.Item(11,3).Value = '南昌·SuperComic动漫游戏博览会'
$ws4.Cells.Item(11,4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Cells.Item(11,5).Value = '2024.07.13 09:00-07.14 17:00'
$ws4.Cells.Item(11,6).Value = 3388
$ws4.Cells.Item(11,7).Value = 65
$ws4.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86992'
$ws4.Cells.Item(11,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/wQTAjelJ1717642148929.jpeg'

$ws4.Cells.Item(12,2).Value = '2024-07-13'
$ws4.Cells.Item(12,3).Value = '南昌·SuperComic配音演员刘明月专场见面会'
$ws4.Cells.Item(12,4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Cells.Item(12,5).Value = '2024.07.13 09:00-07.13 17:00'
$ws4.Cells.Item(12,6).Value = 110
$ws4.Cells.Item(12,7).Value = 168
$ws4.Cells.Item(12,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87570'
$ws4.Cells.Item(12,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/1D1reIl81718609013880.png'

$ws4.Cells.Item(13,2).Value = '2024-07-13'
$ws4.Cells.Item(13,3).Value = '南昌·THO-梦违赣鄱荟萃·叁~幻想Strawberry~!!'
$ws4.Cells.Item(13,4).Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws4.Cells.Item(13,5).Value = '2024.07.13 09:30-07.13 17:30'
$ws4.Cells.Item(13,6).Value = 77
$ws4.Cells.Item(13,7).Value = 65
$ws4.Cells.Item(13,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87668'
$ws4.Cells.Item(13,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/Bk9cYryT1718360290362.jpeg'

$ws4.Cells.Item(14,2).Value = '2024-07-13'
$ws4.Cells.Item(14,3).Value = '宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华'
$ws4.Cells.Item(14,4).Value = '宜春国际商贸城会展中心 宜春国际商贸城会展中心'
$ws4.Cells.Item(14,5).Value = '2024.07.13 10:00-07.14 17:00'
$ws4.Cells.Item(14,6).Value = 63
$ws4.Cells.Item(14,7).Value = 55
$ws4.Cells.Item(14,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86667'
$ws4.Cells.Item(14,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg'

$ws4.Cells.Item(15,2).Value = '2024-07-13'
$ws4.Cells.Item(15,3).Value = '江西·东方LiveParty×THO03幻想Strawberry~！！'
$ws4.Cells.Item(15,4).Value = '上海路543号520Park文创公园21号01区域 瓦肆VAS NANCHANG'
$ws4.Cells.Item(15,5).Value = '2024.07.13 20:30-07.13 23:00'
$ws4.Cells.Item(15,6).Value = 66
$ws4.Cells.Item(15,7).Value = 80
$ws4.Cells.Item(15,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87366'
$ws4.Cells.Item(15,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/1L3I6Qmg1718292516616.jpeg'

$ws4.Cells.Item(16,2).Value = '2024-07-13'
$ws4.Cells.Item(16,3).Value = '赣州·十万伏特-次元交流会（夏）'
$ws4.Cells.Item(16,4).Value = '梅关大道36-16号 麋鹿星球艺术中心'
$ws4.Cells.Item(16,5).Value = '2024.07.13 09:30-07.13 17:00'
$ws4.Cells.Item(16,6).Value = 36
$ws4.Cells.Item(16,7).Value = 45
$ws4.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87597'
$ws4.Cells.Item(16,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/87yQ4Hmf1718681348727.jpeg'

$ws4.Cells.Item(17,2).Value = '2024-07-14'
$ws4.Cells.Item(17,3).Value = '南昌·赛马娘ONLY'
$ws4.Cells.Item(17,4).Value = '洪城路99号 锦都皇冠酒店(八一广场火车站店)'
$ws4.Cells.Item(17,5).Value = '2024.07.14 09:00-07.14 17:30'
$ws4.Cells.Item(17,6).Value = 50
$ws4.Cells.Item(17,7).Value = 68
$ws4.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87367'
$ws4.Cells.Item(17,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/wXQuIKtu1718165450704.png'

$ws4.Cells.Item(18,2).Value = '2024-07-14'
$ws4.Cells.Item(18,3).Value = '吉安·COMIC LIFE次元假日05'
$ws4.Cells.Item(18,4).Value = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws4.Cells.Item(18,5).Value = '2024.07.14 09:00-07.14 18:00'
$ws4.Cells.Item(18,6).Value = 577
$ws4.Cells.Item(18,7).Value = 52.1
$ws4.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85924'
$ws4.Cells.Item(18,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg'

$ws4.Cells.Item(19,2).Value = '2024-07-19'
$ws4.Cells.Item(19,3).Value = '九江·第一届Loading加载中动漫展'
$ws4.Cells.Item(19,4).Value = '湓浦街道大中路339号 百嘉洲际酒店'
$ws4.Cells.Item(19,5).Value = '2024.07.19 09:00-07.21 17:00'
$ws4.Cells.Item(19,6).Value = 62
$ws4.Cells.Item(19,7).Value = 36.6
$ws4.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87787'
$ws4.Cells.Item(19,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/TH3lVD5G1718158901239.jpeg'

$ws4.Cells.Item(20,2).Value = '2024-07-19'
$ws4.Cells.Item(20,3).Value = '赣州·第四届赣州半夏动漫展'
$ws4.Cells.Item(20,4).Value = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws4.Cells.Item(20,5).Value = '2024.07.19 10:00-07.21 17:00'
$ws4.Cells.Item(20,6).Value = 686
$ws4.Cells.Item(20,7).Value = 55
$ws4.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86587'
$ws4.Cells.Item(20,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg'

$ws4.Cells.Item(21,2).Value = '2024-07-20'
$ws4.Cells.Item(21,3).Value = '南昌·漫拥动漫嘉年华Pro-追光启航'
$ws4.Cells.Item(21,4).Value = '小蓝南路420号 洪州体育馆'
$ws4.Cells.Item(21,5).Value = '2024.07.20 09:00-07.21 17:00'
$ws4.Cells.Item(21,6).Value = 198
$ws4.Cells.Item(21,7).Value = 52.5
$ws4.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85796'
$ws4.Cells.Item(21,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png'

$ws4.Cells.Item(22,2).Value = '2024-07-21'
$ws4.Cells.Item(22,3).Value = '乐平·CY境界次元动漫夏时庆'
$ws4.Cells.Item(22,4).Value = '翥山西路182号 佳佳基大酒店'
$ws4.Cells.Item(22,5).Value = '2024.07.21 10:00-07.21 17:00'
$ws4.Cells.Item(22,6).Value = 111
$ws4.Cells.Item(22,7).Value = 35
$ws4.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86768'
$ws4.Cells.Item(22,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png'

$ws4.Cells.Item(23,2).Value = '2024-07-21'
$ws4.Cells.Item(23,3).Value = '九江·SXD动漫嘉年华'
$ws4.Cells.Item(23,4).Value = '湓浦街道大中路339号 百嘉洲际酒店'
$ws4.Cells.Item(23,5).Value = '2024.07.21 10:00-07.21 17:30'
$ws4.Cells.Item(23,6).Value = 54
$ws4.Cells.Item(23,7).Value = 45
$ws4.Cells.Item(23,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86832'
$ws4.Cells.Item(23,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg'

$ws4.Cells.Item(24,2).Value = '2024-07-21'
$ws4.Cells.Item(24,3).Value = '抚州·临次元08·盛夏动漫狂欢节'
$ws4.Cells.Item(24,4).Value = '伍塘路1098号 乐课篮球公园'
$ws4.Cells.Item(24,5).Value = '2024.07.21 10:00-07.21 16:00'
$ws4.Cells.Item(24,6).Value = 50
$ws4.Cells.Item(24,7).Value = 50
$ws4.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87763'
$ws4.Cells.Item(24,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/6qgetbCh1718720523395.jpeg'

$ws4.Cells.Item(25,2).Value = '2024-07-21'
$ws4.Cells.Item(25,3).Value = '萍乡·NL14动漫游戏展·夏日狂想曲'
$ws4.Cells.Item(25,4).Value = '公园南路168号(近工行城北分理处) 梅生嘉华酒店'
$ws4.Cells.Item(25,5).Value = '2024.07.21 10:00-07.21 17:00'
$ws4.Cells.Item(25,6).Value = 59
$ws4.Cells.Item(25,7).Value = 40
$ws4.Cells.Item(25,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86658'
$ws4.Cells.Item(25,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg'

$ws4.Cells.Item(26,2).Value = '2024-07-26'
$ws4.Cells.Item(26,3).Value = '南昌·萌卡动漫展'
$ws4.Cells.Item(26,4).Value = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$ws4.Cells.Item(26,5).Value = '2024.07.26 09:00-07.28 17:00'
$ws4.Cells.Item(26,6).Value = 2447
$ws4.Cells.Item(26,7).Value = 58.5
$ws4.Cells.Item(26,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86776'
$ws4.Cells.Item(26,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg'

$ws4.Cells.Item(27,2).Value = '2024-07-27'
$ws4.Cells.Item(27,3).Value = '江西·次元星河动漫游戏嘉年华'
$ws4.Cells.Item(27,4).Value = '九龙大道1177号 南昌绿地国际博览中心'
$ws4.Cells.Item(27,5).Value = '2024.07.27 10:00-07.28 17:00'
$ws4.Cells.Item(27,6).Value = 4969
$ws4.Cells.Item(27,7).Value = 69
$ws4.Cells.Item(27,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85493'
$ws4.Cells.Item(27,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png'

$ws4.Cells.Item(28,2).Value = '2024-07-27'
$ws4.Cells.Item(28,3).Value = '赣州·马娘only'
$ws4.Cells.Item(28,4).Value = '火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)'
$ws4.Cells.Item(28,5).Value = '2024.07.27 09:00-07.27 17:00'
$ws4.Cells.Item(28,6).Value = 32
$ws4.Cells.Item(28,7).Value = 60
$ws4.Cells.Item(28,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86772'
$ws4.Cells.Item(28,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png'

$ws4.Cells.Item(29,2).Value = '2024-07-28'
$ws4.Cells.Item(29,3).Value = '赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会'
$ws4.Cells.Item(29,4).Value = '兴国路恒大帝景西门 江西长庚控股有限公司'
$ws4.Cells.Item(29,5).Value = '2024.07.28 11:00-07.28 17:00'
$ws4.Cells.Item(29,6).Value = 69
$ws4.Cells.Item(29,7).Value = 56
$ws4.Cells.Item(29,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85688'
$ws4.Cells.Item(29,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png'

$ws4.Cells.Item(30,2).Value = '2024-07-30'
$ws4.Cells.Item(30,3).Value = '宜春·第三十五届静卿国风动漫文化展览会'
$ws4.Cells.Item(30,4).Value = '宜阳大道19号(交通银行旁) 宜春安缦文华酒店'
$ws4.Cells.Item(30,5).Value = '2024.07.30 09:00-07.30 17:00'
$ws4.Cells.Item(30,6).Value = 472
$ws4.Cells.Item(30,7).Value = 45
$ws4.Cells.Item(30,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86684'
$ws4.Cells.Item(30,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/45bGPXfQ1716709212619.jpeg'

$ws4.Cells.Item(31,2).Value = '2024-08-03'
$ws4.Cells.Item(31,3).Value = '南昌·幻梦境国际动漫游戏嘉年华1th'
$ws4.Cells.Item(31,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws4.Cells.Item(31,5).Value = '2024.08.03 09:00-08.04 17:30'
$ws4.Cells.Item(31,6).Value = 1277
$ws4.Cells.Item(31,7).Value = 64
$ws4.Cells.Item(31,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83980'
$ws4.Cells.Item(31,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'

$ws4.Cells.Item(32,2).Value = '2024-08-03'
$ws4.Cells.Item(32,3).Value = '吉安·COMIC LIFE周年庆典'
$ws4.Cells.Item(32,4).Value = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$ws4.Cells.Item(32,5).Value = '2024.08.03 09:30-08.03 18:00'
$ws4.Cells.Item(32,6).Value = 277
$ws4.Cells.Item(32,7).Value = 46.6
$ws4.Cells.Item(32,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87164'
$ws4.Cells.Item(32,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg'

$ws4.Cells.Item(33,2).Value = '2024-08-03'
$ws4.Cells.Item(33,3).Value = '景德镇·第十五届瓷都ACG动漫游戏博览会'
$ws4.Cells.Item(33,4).Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws4.Cells.Item(33,5).Value = '2024.08.03 09:00-08.04 17:00'
$ws4.Cells.Item(33,6).Value = 2199
$ws4.Cells.Item(33,7).Value = 55
$ws4.Cells.Item(33,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86341'
$ws4.Cells.Item(33,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png'

$ws4.Cells.Item(34,2).Value = '2024-08-03'
$ws4.Cells.Item(34,3).Value = '景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票'
$ws4.Cells.Item(34,4).Value = '迎宾大道与寺山路交叉口东200米 陶博城'
$ws4.Cells.Item(34,5).Value = '2024.08.03 08:30-08.03 17:00'
$ws4.Cells.Item(34,6).Value = 570
$ws4.Cells.Item(34,7).Value = '已售罄'
$ws4.Cells.Item(34,8).Value = 'https://show.bilibili.com/platform/detail.html?id=85981'
$ws4.Cells.Item(34,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png'

$ws4.Cells.Item(35,2).Value = '2024-08-03'
$ws4.Cells.Item(35,3).Value = '樟树·第二届静卿国风动漫文化展览会'
$ws4.Cells.Item(35,4).Value = '杏佛路89号 樟树银河国际酒店'
$ws4.Cells.Item(35,5).Value = '2024.08.03 09:00-08.03 17:00'
$ws4.Cells.Item(35,6).Value = 484
$ws4.Cells.Item(35,7).Value = 45
$ws4.Cells.Item(35,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86683'
$ws4.Cells.Item(35,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg'

$ws4.Cells.Item(36,2).Value = '2024-08-03'
$ws4.Cells.Item(36,3).Value = '萍乡·AU9夏至国漫展'
$ws4.Cells.Item(36,4).Value = '安源中大道17号 壹号公馆（萍乡）'
$ws4.Cells.Item(36,5).Value = '2024.08.03 10:00-08.03 17:00'
$ws4.Cells.Item(36,6).Value = 75
$ws4.Cells.Item(36,7).Value = 45
$ws4.Cells.Item(36,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86453'
$ws4.Cells.Item(36,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/hm1EACno1718936156944.jpeg'

$ws4.Cells.Item(37,2).Value = '2024-08-03'
$ws4.Cells.Item(37,3).Value = '赣州·第一届环梦动漫游戏嘉年华'
$ws4.Cells.Item(37,4).Value = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$ws4.Cells.Item(37,5).Value = '2024.08.03 09:00-08.05 17:00'
$ws4.Cells.Item(37,6).Value = 85
$ws4.Cells.Item(37,7).Value = 36.6
$ws4.Cells.Item(37,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87449'
$ws4.Cells.Item(37,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/iC3PuUfR1717740188790.jpeg'

$ws4.Cells.Item(38,2).Value = '2024-08-04'
$ws4.Cells.Item(38,3).Value = '上饶·第十五届IX Group国风嘉年华暨十周年庆典'
$ws4.Cells.Item(38,4).Value = '高铁经济试验区凤凰东大道16号7幢 上饶饶商金茂诚悦酒店(上饶高铁站)'
$ws4.Cells.Item(38,5).Value = '2024.08.04 09:30-08.04 17:30'
$ws4.Cells.Item(38,6).Value = 156
$ws4.Cells.Item(38,7).Value = 60
$ws4.Cells.Item(38,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87225'
$ws4.Cells.Item(38,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/l5fIXZSX1717562269098.jpeg'

$ws4.Cells.Item(39,2).Value = '2024-08-04'
$ws4.Cells.Item(39,3).Value = '九江·第一届异次元动漫嘉年华'
$ws4.Cells.Item(39,4).Value = '长虹西大道兴城广场99号 九江半岛宾馆'
$ws4.Cells.Item(39,5).Value = '2024.08.04 08:00-08.04 17:00'
$ws4.Cells.Item(39,6).Value = 305
$ws4.Cells.Item(39,7).Value = 45
$ws4.Cells.Item(39,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84407'
$ws4.Cells.Item(39,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg'

$ws4.Cells.Item(40,2).Value = '2024-08-06'
$ws4.Cells.Item(40,3).Value = '南昌·第一届异次元动漫嘉年华'
$ws4.Cells.Item(40,4).Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws4.Cells.Item(40,5).Value = '2024.08.06 08:00-08.06 17:00'
$ws4.Cells.Item(40,6).Value = 452
$ws4.Cells.Item(40,7).Value = 55
$ws4.Cells.Item(40,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84102'
$ws4.Cells.Item(40,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg'

$ws4.Cells.Item(41,2).Value = '2024-08-08'
$ws4.Cells.Item(41,3).Value = '赣州·第二届异次元动漫嘉年华'
$ws4.Cells.Item(41,4).Value = '金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆'
$ws4.Cells.Item(41,5).Value = '2024.08.08 08:00-08.08 17:00'
$ws4.Cells.Item(41,6).Value = 777
$ws4.Cells.Item(41,7).Value = 45
$ws4.Cells.Item(41,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84184'
$ws4.Cells.Item(41,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg'

$ws4.Cells.Item(42,2).Value = '2024-08-10'
$ws4.Cells.Item(42,3).Value = '南昌·花绒万兽第二聚'
$ws4.Cells.Item(42,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws4.Cells.Item(42,5).Value = '2024.08.10 10:00-08.11 17:00'
$ws4.Cells.Item(42,6).Value = 25
$ws4.Cells.Item(42,7).Value = 188
$ws4.Cells.Item(42,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87600'
$ws4.Cells.Item(42,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/i0Ojsne01718693886054.png'

$ws4.Cells.Item(43,2).Value = '2024-08-10'
$ws4.Cells.Item(43,3).Value = '高安·第二届静卿国风动漫文化展览会'
$ws4.Cells.Item(43,4).Value = '华林中路606号 高安华鼎国际大酒店'
$ws4.Cells.Item(43,5).Value = '2024.08.10 09:00-08.10 17:00'
$ws4.Cells.Item(43,6).Value = 449
$ws4.Cells.Item(43,7).Value = 45
$ws4.Cells.Item(43,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86682'
$ws4.Cells.Item(43,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg'

$ws4.Cells.Item(44,2).Value = '2024-08-15'
$ws4.Cells.Item(44,3).Value = '上饶·次元重现夏日嘉年华'
$ws4.Cells.Item(44,4).Value = '普济巷地委大院北侧约90米 四季体育运动馆'
$ws4.Cells.Item(44,5).Value = '2024.08.15 09:30-08.15 17:30'
$ws4.Cells.Item(44,6).Value = 27
$ws4.Cells.Item(44,7).Value = 48
$ws4.Cells.Item(44,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87679'
$ws4.Cells.Item(44,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/fxlKV2SL1718784421064.jpeg'

$ws4.Cells.Item(45,2).Value = '2024-08-24'
$ws4.Cells.Item(45,3).Value = '南昌·第四届龙年动漫展——暑假最后的狂欢'
$ws4.Cells.Item(45,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws4.Cells.Item(45,5).Value = '2024.08.24 10:00-08.25 18:00'
$ws4.Cells.Item(45,6).Value = 458
$ws4.Cells.Item(45,7).Value = 45
$ws4.Cells.Item(45,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87135'
$ws4.Cells.Item(45,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'

$ws4.Rows.Item(46).Delete()

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,6).Value = 66
